$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update "Riders" (column C) and "Average" (column D) values for the week
# new Madigan bike hours

$ws.Range("C2").Value = 273
$ws.Range("D2").Value = 264

$ws.Range("C3").Value = 218
$ws.Range("D3").Value = 230.83

$ws.Range("C4").Value = 204
$ws.Range("D4").Value = 198.33

$ws.Range("C5").Value = 137
$ws.Range("D5").Value = 226.4

$ws.Range("C6").Value = 80
$ws.Range("D6").Value = 211

$ws.Range("C7").Value = 141
$ws.Range("D7").Value = 125.57

$ws.Range("C8").Value = 162
$ws.Range("D8").Value = 99.86

$wb.Save()
